{"js": "// Lattice multiplication worksheet: swap in a new set of \"A x B\" exercises.\n// Each table cell holds a single paragraph/run whose text (split on line\n// breaks) looks like:\n//   \"A x B\"\n//   \"  {tens(B)}    {ones(B)}\"\n//   \"  ----\"\n//   \"{tens(A)}|    |\"\n//   \"{ones(A)}|    |\"\n// The new (A, B) pairs below are applied in table reading order (row major,\n// left to right, top to bottom) -- the table shape (5 rows x 3 cols) itself\n// does not change, only the text inside each existing cell.\n\nconst newPairs = [\n  [71, 11], [49, 81], [83, 57],\n  [98, 56], [95, 66], [31, 47],\n  [14, 29], [73, 66], [71, 61],\n  [30, 69], [63, 91], [31, 24],\n  [76, 72], [41, 66], [21, 83],\n];\n\nconst VB = String.fromCharCode(11); // vertical-tab char Word stores <w:br/> as\n\nfunction cellText(a, b) {\n  const aTens = String(Math.floor(a / 10));\n  const aOnes = String(a % 10);\n  const bTens = String(Math.floor(b / 10));\n  const bOnes = String(b % 10);\n  return [\n    `${a} x ${b}`,\n    `  ${bTens}    ${bOnes}`,\n    \"  ----\",\n    `${aTens}|    |`,\n    `${aOnes}|    |`,\n  ].join(VB);\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst colCount = table.values[0].length;\n\nfor (let i = 0; i < newPairs.length; i++) {\n  const row = Math.floor(i / colCount);\n  const col = i % colCount;\n  const cell = table.getCell(row, col);\n  const body = cell.body;\n  body.paragraphs.load(\"items\");\n  await context.sync();\n\n  const [a, b] = newPairs[i];\n  const paragraph = body.paragraphs.items[0];\n  paragraph.insertText(cellText(a, b), Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Lattice multiplication worksheet: swap in a new set of \"A x B\" exercises.\n# Each table cell holds a single paragraph whose text (split on line breaks,\n# i.e. Word's vertical-tab char 11 produced by <w:br/>) looks like:\n#   \"A x B\"\n#   \"  {tens(B)}    {ones(B)}\"\n#   \"  ----\"\n#   \"{tens(A)}|    |\"\n#   \"{ones(A)}|    |\"\n# The new (A, B) pairs below are applied in table reading order (row major,\n# left to right, top to bottom) -- the table shape (5 rows x 3 cols) itself\n# does not change, only the text inside each existing cell.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n$colCount = $table.Columns.Count\n\n$newPairs = @(\n    @(71, 11), @(49, 81), @(83, 57),\n    @(98, 56), @(95, 66), @(31, 47),\n    @(14, 29), @(73, 66), @(71, 61),\n    @(30, 69), @(63, 91), @(31, 24),\n    @(76, 72), @(41, 66), @(21, 83)\n)\n\n$vb = [char]11\n\nfor ($i = 0; $i -lt $newPairs.Count; $i++) {\n    $row = [int][Math]::Floor($i / $colCount) + 1\n    $col = ($i % $colCount) + 1\n\n    $a = $newPairs[$i][0]\n    $b = $newPairs[$i][1]\n\n    $aTens = [int][Math]::Floor($a / 10)\n    $aOnes = $a % 10\n    $bTens = [int][Math]::Floor($b / 10)\n    $bOnes = $b % 10\n\n    $text = \"$a x $b\" + $vb + \"  $bTens    $bOnes\" + $vb + \"  ----\" + $vb + \"$aTens|    |\" + $vb + \"$aOnes|    |\"\n\n    $table.Cell($row, $col).Range.Text = $text\n}\n"}
